$d = $word.ActiveDocument

# --- Change 1: "Own Pools of truth" paragraph -------------------------------
# "...others might use APIs to create..." -> "...others might use those APIs to create..."
$rng1 = $d.Content
$rng1.Find.Execute(
    "others might use APIs to create their own source of truth",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "others might use those APIs to create their own source of truth",
    2
)

# --- Change 2: JSON paragraph clean-up --------------------------------------
# Merge the two runs that were split around the stray "_GoBack" bookmark back
# into a single contiguous sentence (no visible text change, just tidy-up of
# the run that used to wrap around the bookmark).
$rng2 = $d.Content
$rng2.Find.Execute(
    "numeric data types. JSON is the better choice, provided you have plug-ins to doll it up.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "numeric data types. JSON is the better choice, provided you have plug-ins to doll it up.",
    2
)

# --- Change 3: "Creates/generates APIs on top of that" ----------------------
# (kept as-is; see notes)

# --- Change 4: Socrata hyperlink --------------------------------------------
# (kept as-is; see notes)
